$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet grows from 52 data rows (A1:E52) to 53 data rows (A1:E53).
# Row 52 (A52, which carries the date number-format style) is copied into the
# new row 53 first so the new date cell inherits the same style (s="2"),
# then every cell A2:E53 is overwritten with the corrected/updated values below.
$ws.Range("A52").Copy($ws.Range("A53"))

$data = @(
    @(2, 39400, 2007, 0.4235526809466261, 2008, 0.6439341879002525),
    @(3, 39583, 2008, -0.3623658873974311, 2009, 0.1825419310453658),
    @(4, 39765, 2008, -0.5718076928962645, 2009, -0.1800933741311961),
    @(5, 39948, 2009, -0.009261555895478946, 2010, 0.1145211022186787),
    @(6, 40130, 2009, 0.3486139762225005, 2010, 0.1555182634501051),
    @(7, 40310, 2010, -1.404263945418582, 2011, -0.807808220045203),
    @(8, 40494, 2010, -0.1384957661262898, 2011, 0.6938817570587785),
    @(9, 40676, 2011, 1.692932643509848, 2012, 0.6262577107155831),
    @(10, 40862, 2011, 1.566479473280147, 2012, 0.9614071719361794),
    @(11, 41044, 2012, 1.020829760720643, 2013, 1.148272834981245),
    @(12, 41228, 2012, 0.7307568962936939, 2013, 1.09290550768979),
    @(13, 41409, 2013, 0.6772121200332215, 2014, 1.258913537332895),
    @(14, 41592, 2013, 0.818818812164257, 2014, 0.9607602172681418),
    @(15, 41774, 2014, 1.019715257608911, 2015, 0.9536145745415947),
    @(16, 41957, 2014, 0.9180054319587239, 2015, 1.375398114243209),
    @(17, 42137, 2015, 2.173959184500385, 2016, 1.566646323486043),
    @(18, 42321, 2015, 1.984684278296656, 2016, 1.473274087935805),
    @(19, 42503, 2016, 1.707434489469994, 2017, 1.30258347990615),
    @(20, 42689, 2016, 1.755995812646982, 2017, 1.681032827388362),
    @(21, 42867, 2017, 1.456988786619839, 2018, 1.842797144428188),
    @(22, 43053, 2017, 1.946965557828384, 2018, 1.755491062323111),
    @(23, 43145, 2018, 1.131202984360957, 2019, 1.657737120813452),
    @(24, 43235, 2018, 1.241332692055597, 2019, 1.58004210678635),
    @(25, 43326, 2018, 1.260396653238383, 2019, 1.567743002885069),
    @(26, 43418, 2018, 1.06432145354225, 2019, 0.776718238020746),
    @(27, 43510, 2019, 0.5757500748109434, 2020, 1.030688008679626),
    @(28, 43600, 2019, 1.592885137608979, 2020, 1.604795846351514),
    @(29, 43691, 2019, 1.308235387832934, 2020, 1.242807488305719),
    @(30, 43783, 2019, 1.361817904277696, 2020, 1.316199564471554),
    @(31, 43875, 2020, 0.9437384066259158, 2021, 0.904959070968947),
    @(32, 43966, 2020, -2.015335584265165, 2021, -1.215549235925828),
    @(33, 44068, 2020, -5.210209911466245, 2021, -2.349089443609143),
    @(34, 44159, 2020, -4.352425014431304, 2021, 0.03547044462246518),
    @(35, 44251, 2021, -4.454337270215236, 2022, -3.012953608516933),
    @(36, 44341, 2021, -3.579597300369253, 2022, -1.403103901755631),
    @(37, 44432, 2021, -1.897775264882628, 2022, 4.997412520017441),
    @(38, 44525, 2021, -1.761645650979182, 2022, 3.765721202592909),
    @(39, 44617, 2022, 3.304925622412869, 2023, -0.4653479251390702),
    @(40, 44706, 2022, 4.461954539041502, 2023, 0.7797949948739058),
    @(41, 44798, 2022, 4.787836378515364, 2023, 1.112953228460167),
    @(42, 44890, 2022, 5.20787683103745, 2023, 3.217995704408838),
    @(43, 44981, 2023, -0.5032298616872488, 2024, 0.34496349151385),
    @(44, 45071, 2023, -1.305206755692701, 2024, 0.5821000732047832),
    @(45, 45163, 2023, -0.49899188013105, 2024, 2.610266500707703),
    @(46, 45254, 2023, -0.9008525709169546, 2024, 0.6027009207580036),
    @(47, 45345, 2024, 0.2229020320597241, 2025, -0.08457672677967265),
    @(48, 45436, 2024, 0.0845726262934221, 2025, 0.1341520870597357),
    @(49, 45534, 2024, 0.5084754301873051, 2025, 0.0148008406940292),
    @(50, 45618, 2024, 0.2738544794132824, 2025, 0.2681899963140832),
    @(51, 45713, 2025, 0.2312068876759277, 2026, -0.4074819591325718),
    @(52, 45800, 2025, 0.9724700385226326, 2026, 0.6236501628417823),
    @(53, 45891, 2025, 0.9584581489103794, 2026, 0.6809779381435677)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
}
